$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "23.13") are
# stored as text, matching the original inlineStr cell type, not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{Row=2; D="27.905.99"; E="  -0.07%  "},
    @{Row=3; D="1.626.48"; E="  -0.26%  "},
    @{Row=4; E="  +0.11%  "},
    @{Row=5; D="211.58"; E="  -0.17%  "},
    @{Row=6; D="0.516"; E="  -1.29%  "},
    @{Row=7; E="  +0.11%  "},
    @{Row=8; D="23.13"; E="  -0.90%  "},
    @{Row=9; E="  +0.30%  "},
    @{Row=10; D="0.0606"; E="  -1.25%  "},
    @{Row=11; E="  +0.23%  "},
    @{Row=12; D="1.859.70"; E="  -0.12%  "},
    @{Row=13; D="1.627.57"; E="  -0.02%  "},
    @{Row=15; E="  -1.37%  "},
    @{Row=16; D="64.70"; E="  -1.37%  "},
    @{Row=17; D="27.922.12"; E="  +0.01%  "},
    @{Row=18; D="227.69"; E="  -1.11%  "},
    @{Row=19; D="7.61"; E="  -0.45%  "},
    @{Row=20; D="0.0₃0716"; E="  -0.97%  "},
    @{Row=21; E="  +0.10%  "},
    @{Row=22; E="  -0.10%  "},
    @{Row=23; D="9.95"; E="  -3.09%  "},
    @{Row=24; E="  +1.34%  "},
    @{Row=25; D="154.28"; E="  -0.25%  "},
    @{Row=26; E="  -0.41%  "},
    @{Row=27; E="  -0.31%  "},
    @{Row=28; E="  +0.23%  "},
    @{Row=29; D="15.38"; E="  -1.12%  "},
    @{Row=30; E="  -0.15%  "},
    @{Row=31; D="0.0480"; E="  -0.35%  "},
    @{Row=32; E="  +0.03%  "},
    @{Row=33; D="1.416.45"; E="  +1.22%  "},
    @{Row=34; E="  +0.66%  "},
    @{Row=35; D="1.61"; E="  +2.02%  "},
    @{Row=36; D="0.999"; E="  -2.48%  "},
    @{Row=37; E="  -0.46%  "},
    @{Row=38; E="  -0.75%  "},
    @{Row=39; D="0.554"; E="  -0.15%  "},
    @{Row=40; D="0.850"; E="  -1.62%  "},
    @{Row=41; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="1.00"; E="  +0.10%  "},
    @{Row=42; B="WEMIXToken"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="1.01"; E="  -1.92%  "},
    @{Row=43; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="1.82"; E="  -0.78%  "},
    @{Row=44; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="65.50"; E="  -1.34%  "},
    @{Row=45; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="5.40"; E="  -1.50%  "},
    @{Row=46; B="RocketPoolETH"; C="https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; D="1.768.78"; E="  -0.17%  "},
    @{Row=47; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.11"; E="  -3.74%  "},
    @{Row=48; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="88.76"; E="  +0.76%  "},
    @{Row=49; B="BabyDogeCoin"; C="https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; D="0.0⁦0102"; E="  -0.50%  "},
    @{Row=50; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.101"; E="  +0.26%  "},
    @{Row=51; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.0503"; E="  -0.47%  "}
)

foreach ($item in $updates) {
    if ($item.ContainsKey("B")) { $ws.Cells.Item($item.Row, 2).Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Cells.Item($item.Row, 3).Value = $item.C }
    if ($item.ContainsKey("D")) { $ws.Cells.Item($item.Row, 4).Value = $item.D }
    if ($item.ContainsKey("E")) { $ws.Cells.Item($item.Row, 5).Value = $item.E }
}

# Clear the explicit number-format style added above so cells end up with the
# same (unstyled) appearance as the rest of the sheet.
$ws.Range("D2:D51").Style = "Normal"
